# Fruta / hortaliza, semanal
# Insert a new weekly record at row 210 (pushing the existing rows 210..231
# down to 211..232) for "Vega Modelo de Temuco" / Granada.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 210; Excel shifts rows 210..231 down to 211..232
# and carries the row's existing formatting (e.g. the date style on column D).
$ws.Rows.Item(210).Insert()

$row = 210

$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 45077
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100104
$ws.Cells.Item($row, 8).Value = "Frutos de pepita"
$ws.Cells.Item($row, 9).Value = 100104001
$ws.Cells.Item($row, 10).Value = "Granada"
$ws.Cells.Item($row, 11).Value = "Wonderfull"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 80
$ws.Cells.Item($row, 14).Value = 22000
$ws.Cells.Item($row, 15).Value = 22000
$ws.Cells.Item($row, 16).Value = 22000
$ws.Cells.Item($row, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 19).Value = 1467
$ws.Cells.Item($row, 20).Value = 15
